# Automatische test-sync: 2025-06-24 21:15:50
# Adds a new logged e-mail (row 28) to the "Logs" sheet, extends the
# conditional formatting ranges to cover it, and bumps the "Retour /
# Terugbetaling" tally on the "Dashboard" sheet from 10 to 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A28").Value = "Omruilen van verkeerd formaat"
$ws.Range("B28").Value = "mailmind.test@zohomail.eu"
$ws.Range("C28").Value = "Hallo,`r`n Ik heb een product ontvangen dat niet de juiste maat is. Is het mogelijk om deze om te ruilen voor de juiste maat? Alvast bedankt!`r`nSent using {0}"
$ws.Range("D28").Value = "Retour / Terugbetaling"
$ws.Range("E28").Value = "Beste klant,`r`nBedankt voor je bericht. We vinden het vervelend om te horen dat je product niet de juiste maat heeft. We helpen je graag met het omruilen van het product voor de juiste maat.`r`nOm het omruilproces te starten, hebben we wat meer informatie nodig. Kun je ons alsjeblieft de volgende gegevens sturen:`r`n- Je ordernummer`r`n- De naam of beschrijving van het product dat je wilt omruilen`r`n- De juiste maat van het product`r`nZodra we deze informatie hebben ontvangen, zullen we zo snel mogelijk contact met je opnemen om de omruiling in gang te zetten.`r`nMet vriendelijke groet,`r`n[E-mailassistent] van [Bedrijfsnaam]"
$ws.Range("F28").Value = "2025-06-24 21:15:08"
$ws.Range("G28").Value = "Ja"

# Extend the two conditional-formatting blocks (Categorie / Beantwoord)
# so they keep covering the full data range through row 28.
$ws.Range("D2:D27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D28"))
$ws.Range("G2:G27").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G28"))

# Update the Dashboard summary count for "Retour / Terugbetaling".
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 11
